# This workbook is a weekly-refreshed price report ("Fruta / hortaliza, semanal").
# The new publication cycle adds 3 fresh daily price entries for the
# "Terminal La Palmera de La Serena - Cereza" sheet. These are inserted as
# new rows 181-183, which pushes the previously-existing rows 181-260 down
# to rows 184-263 (their own contents are unchanged by the insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 181, shifting existing rows 181:260 down to 184:263.
$ws.Rows("181:183").Insert()

# Helper to populate one full data row (columns A:T) in one go.
function Set-CerezaRow {
    param(
        [int]$Row,
        [int]$A, [string]$B, [string]$C, [double]$D, [int]$E, [string]$F,
        [int]$G, [string]$H, [long]$I, [string]$J, [string]$K, [string]$L,
        [int]$M, [int]$N, [int]$O, [int]$P, [string]$Q, [string]$R,
        [int]$S, [int]$T
    )
    $ws.Cells.Item($Row, 1).Value  = $A
    $ws.Cells.Item($Row, 2).Value  = $B
    $ws.Cells.Item($Row, 3).Value  = $C
    $ws.Cells.Item($Row, 4).Value2 = $D
    $ws.Cells.Item($Row, 5).Value  = $E
    $ws.Cells.Item($Row, 6).Value  = $F
    $ws.Cells.Item($Row, 7).Value  = $G
    $ws.Cells.Item($Row, 8).Value  = $H
    $ws.Cells.Item($Row, 9).Value  = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
}

# New row 181: Cereza / Lapins / Primera, Provincia de Curicó, $/bandeja 10 kilos
Set-CerezaRow 181 8 "Terminal La Palmera de La Serena" "Coquimbo" 44572 4 `
    "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Lapins" "Primera" 400 8500 9000 8750 `
    "`$/bandeja 10 kilos" "Provincia de Curicó" 875 10

# New row 182: Cereza / Lapins / Segunda, Provincia de Curicó, $/bandeja 10 kilos
Set-CerezaRow 182 8 "Terminal La Palmera de La Serena" "Coquimbo" 44572 4 `
    "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Lapins" "Segunda" 300 6500 7000 6750 `
    "`$/bandeja 10 kilos" "Provincia de Curicó" 675 10

# New row 183: Cereza / Santina / Primera, Provincia de Curicó, $/caja 15 kilos
Set-CerezaRow 183 8 "Terminal La Palmera de La Serena" "Coquimbo" 44572 4 `
    "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Santina" "Primera" 400 11000 12000 11500 `
    "`$/caja 15 kilos" "Provincia de Curicó" 767 15
